# processa_dados: adiciona a função separa_namostra
# Adds a new "A415" column (E) to the Results sheet, computed as A405 (col D) - 0.023,
# for every data row (3..98), plus the associated header/formatting and sheet-selection
# housekeeping that came along with the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$hdr = $wb.Worksheets.Item("Header")

# --- New header cell: E2 = "A415" (bold, right aligned, like the other headers) ---
$ws.Cells.Item(2, 5).Value = "A415"

# --- New formula column: E3:E98 = D-0.023 ---
# Written as three separate fill passes (E3 alone, then E4:E67, then E68:E98) so the
# resulting shared-formula groups match how the column was built by hand.
$ws.Range("E3").Formula = "=D3-0.023"
$ws.Range("E4:E67").Formula = "=D4-0.023"
$ws.Range("E68:E98").Formula = "=D68-0.023"

# --- Formatting for the new column ---
$ws.Range("E3:E98").NumberFormat = "#####0.000"
$ws.Range("E3:E98").Font.Name = "Arial"
$ws.Range("E3:E98").Font.Size = 10
$ws.Range("E3:E98").HorizontalAlignment = -4152

$ws.Cells.Item(2, 5).Font.Bold = $true
$ws.Cells.Item(2, 5).HorizontalAlignment = -4152

# --- Selection / active-sheet housekeeping (matches the saved view state) ---
$hdr.Activate()
$hdr.Range("A2").Select()

$ws.Activate()
$ws.Range("G2").Select()
